$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 7573
$ws.Range("J3").Value = 7975
$ws.Range("D4").Value = 1962
$ws.Range("H4").Value = 1711
$ws.Range("J4").Value = 1735
$ws.Range("J5").Value = 621
$ws.Range("J6").Value = 10885
$ws.Range("D7").Value = 28152
$ws.Range("H7").Value = 26022
$ws.Range("J7").Value = 28789

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 228
$ws.Range("J7").Value = 821
$ws.Range("J8").Value = 1817
$ws.Range("J10").Value = 205
$ws.Range("J14").Value = 151
$ws.Range("J18").Value = 231
$ws.Range("J19").Value = 838
$ws.Range("J20").Value = 624
$ws.Range("J21").Value = 83
$ws.Range("J24").Value = 104
$ws.Range("J25").Value = 151
$ws.Range("J27").Value = 175
$ws.Range("D29").Value = 1726
$ws.Range("J29").Value = 1536
$ws.Range("J31").Value = 303
$ws.Range("J33").Value = 1301
$ws.Range("J37").Value = 885
$ws.Range("J41").Value = 221
$ws.Range("J42").Value = 1220
$ws.Range("J44").Value = 225
$ws.Range("J49").Value = 174
$ws.Range("J50").Value = 179
$ws.Range("J51").Value = 361
$ws.Range("J52").Value = 733
$ws.Range("J53").Value = 430
$ws.Range("J54").Value = 562
$ws.Range("J57").Value = 138
$ws.Range("J60").Value = 170
$ws.Range("H63").Value = 269
$ws.Range("J63").Value = 87
$ws.Range("J65").Value = 722
$ws.Range("J67").Value = 1045
$ws.Range("J69").Value = 57
$ws.Range("J72").Value = 112
$ws.Range("J73").Value = 285
$ws.Range("J74").Value = 33
$ws.Range("J78").Value = 337
$ws.Range("J79").Value = 790
$ws.Range("J83").Value = 583
$ws.Range("J85").Value = 1180
$ws.Range("J88").Value = 303
$ws.Range("J89").Value = 359
$ws.Range("J94").Value = 318
$ws.Range("J95").Value = 410
$ws.Range("J98").Value = 212
$ws.Range("D101").Value = 28152
$ws.Range("H101").Value = 26022
$ws.Range("J101").Value = 28789

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J2").Value = 49
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 248
$ws.Range("J6").Value = 260
$ws.Range("J7").Value = 821

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 107
$ws.Range("J7").Value = 359

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 314
$ws.Range("J7").Value = 1180

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 205
$ws.Range("J6").Value = 315
$ws.Range("J7").Value = 733

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 278
$ws.Range("J7").Value = 430

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 481
$ws.Range("J4").Value = 93
$ws.Range("J5").Value = 47
$ws.Range("J6").Value = 676
$ws.Range("J7").Value = 1817

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 176
$ws.Range("J7").Value = 583

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 432
$ws.Range("J6").Value = 465
$ws.Range("J7").Value = 1301

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 148
$ws.Range("J7").Value = 410

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 268
$ws.Range("J3").Value = 298
$ws.Range("J7").Value = 885

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 271
$ws.Range("J7").Value = 722

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 98
$ws.Range("J7").Value = 303

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J6").Value = 290
$ws.Range("J7").Value = 1045

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 174

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 142
$ws.Range("J7").Value = 562

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 465
$ws.Range("J3").Value = 541
$ws.Range("D4").Value = 79
$ws.Range("J6").Value = 391
$ws.Range("D7").Value = 1726
$ws.Range("J7").Value = 1536

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J4").Value = 42
$ws.Range("J7").Value = 838

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 255
$ws.Range("J3").Value = 246
$ws.Range("J6").Value = 645
$ws.Range("J7").Value = 1220

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 118
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 104
$ws.Range("J6").Value = 105
$ws.Range("J7").Value = 337

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 228
$ws.Range("J7").Value = 790

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 205
$ws.Range("J7").Value = 624

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 61
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 127
$ws.Range("J4").Value = 17

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 170
$ws.Range("J7").Value = 318

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 64
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 132
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J2").Value = 47
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 285

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 67
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 161
$ws.Range("J7").Value = 303

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J6").Value = 149
$ws.Range("J7").Value = 361

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 138

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 58
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 33
